$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting the existing row 31 (and below) down to row 32.
$ws.Rows.Item(31).Insert()

# Copy formatting for the new row 31 from row 30 (so column D keeps its date style, etc.)
$ws.Range("A30:T30").Copy()
$ws.Range("A31:T31").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 31 with the values that row 30 held before this edit
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = (Get-Date -Year 2023 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100104
$ws.Range("H31").Value = "Frutos de pepita"
$ws.Range("I31").Value = 100104003
$ws.Range("J31").Value = "Membrillo"
$ws.Range("K31").Value = "Champion"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 9000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 9500
$ws.Range("Q31").Value = '$/bandeja 18 kilos granel'
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 528
$ws.Range("T31").Value = 18

# Update row 30 with its new values
$ws.Range("D30").Value = (Get-Date -Year 2023 -Month 8 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N30").Value = 13000
$ws.Range("O30").Value = 14000
$ws.Range("P30").Value = 13500
$ws.Range("S30").Value = 750

$wb.Save()
